# Applies numeric updates to the Leve profit-tracking sheets
# (currentAveragePrice* / LevePrice* / LeveProfit* columns), refreshed
# from the latest market-board snapshot.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 218069.17
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 218069.17
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 654207.51
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -654543.51
$ws.Range("H33").Value = 721.9
$ws.Range("I33").Value = 396.6154
$ws.Range("J33").Value = 1326
$ws.Range("K33").Value = 396.6154
$ws.Range("L33").Value = 1326
$ws.Range("M33").Value = -167.6154
$ws.Range("N33").Value = -1784
$ws.Range("H100").Value = 1952.091
$ws.Range("I100").Value = 1252.6666
$ws.Range("J100").Value = 5099.5
$ws.Range("K100").Value = 1252.6666
$ws.Range("L100").Value = 5099.5
$ws.Range("M100").Value = -711.6666
$ws.Range("N100").Value = -6181.5
$ws.Range("H112").Value = 60316.65
$ws.Range("I112").Value = 334120.34
$ws.Range("K112").Value = 1002361.02
$ws.Range("M112").Value = -1001253.02
$ws.Range("H113").Value = 333338020
$ws.Range("I113").Value = 500002500
$ws.Range("J113").Value = 9000
$ws.Range("K113").Value = 500002500
$ws.Range("L113").Value = 9000
$ws.Range("M113").Value = -499999246
$ws.Range("N113").Value = -15508
$ws.Range("H131").Value = 5616.8184
$ws.Range("I131").Value = 976.1111
$ws.Range("K131").Value = 2928.3333
$ws.Range("M131").Value = 2111.6667
$ws.Range("H132").Value = 3353.0278
$ws.Range("I132").Value = 3345.2942
$ws.Range("K132").Value = 10035.8826
$ws.Range("M132").Value = -7505.882599999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3779.3389
$ws.Range("I32").Value = 2195.7646
$ws.Range("K32").Value = 2195.7646
$ws.Range("M32").Value = -1908.7646
$ws.Range("H45").Value = 1396.7778
$ws.Range("I45").Value = 667.2857
$ws.Range("K45").Value = 667.2857
$ws.Range("M45").Value = -290.2857
$ws.Range("H61").Value = 41667520
$ws.Range("I61").Value = 41667520
$ws.Range("K61").Value = 41667520
$ws.Range("M61").Value = -41667308
$ws.Range("H74").Value = 27779660
$ws.Range("I74").Value = 29413524
$ws.Range("K74").Value = 29413524
$ws.Range("M74").Value = -29412650
$ws.Range("H77").Value = 27779660
$ws.Range("I77").Value = 29413524
$ws.Range("K77").Value = 147067620
$ws.Range("M77").Value = -147063252
$ws.Range("H110").Value = 111114970
$ws.Range("I110").Value = 125003720
$ws.Range("K110").Value = 125003720
$ws.Range("M110").Value = -125001675
$ws.Range("H122").Value = 14495213
$ws.Range("I122").Value = 19609588
$ws.Range("K122").Value = 58828764
$ws.Range("M122").Value = -58826314
$ws.Range("H132").Value = 1891.3
$ws.Range("I132").Value = 1620.375
$ws.Range("J132").Value = 2975
$ws.Range("K132").Value = 4861.125
$ws.Range("L132").Value = 8925
$ws.Range("M132").Value = -2331.125
$ws.Range("N132").Value = -13985
$ws.Range("H136").Value = 41667520
$ws.Range("I136").Value = 41667520
$ws.Range("K136").Value = 125002560
$ws.Range("M136").Value = -125000010

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 29415928
$ws.Range("I86").Value = 41671252
$ws.Range("J86").Value = 3149.4
$ws.Range("K86").Value = 41671252
$ws.Range("L86").Value = 3149.4
$ws.Range("M86").Value = -41670129
$ws.Range("N86").Value = -5395.4
$ws.Range("H89").Value = 29415928
$ws.Range("I89").Value = 41671252
$ws.Range("J89").Value = 3149.4
$ws.Range("K89").Value = 208356260
$ws.Range("L89").Value = 15747
$ws.Range("M89").Value = -208350644
$ws.Range("N89").Value = -26979
$ws.Range("H99").Value = 2503.25
$ws.Range("I99").Value = 2079.75
$ws.Range("K99").Value = 2079.75
$ws.Range("M99").Value = -581.75
$ws.Range("H107").Value = 166667870
$ws.Range("I107").Value = 1809.5
$ws.Range("J107").Value = 500000000
$ws.Range("K107").Value = 1809.5
$ws.Range("L107").Value = 500000000
$ws.Range("M107").Value = 110.5
$ws.Range("N107").Value = -500003840
$ws.Range("H132").Value = 94393
$ws.Range("J132").Value = 94393
$ws.Range("L132").Value = 94393
$ws.Range("N132").Value = -104513

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2606.7273
$ws.Range("I16").Value = 1794.8
$ws.Range("K16").Value = 1794.8
$ws.Range("M16").Value = -1507.8
$ws.Range("H31").Value = 3211.1936
$ws.Range("J31").Value = 3564.647
$ws.Range("L31").Value = 3564.647
$ws.Range("N31").Value = -4154.647
$ws.Range("H34").Value = 3211.1936
$ws.Range("J34").Value = 3564.647
$ws.Range("L34").Value = 3564.647
$ws.Range("N34").Value = -3968.647
$ws.Range("H62").Value = 9332.556
$ws.Range("I62").Value = 7999.25
$ws.Range("K62").Value = 7999.25
$ws.Range("M62").Value = -7375.25
$ws.Range("H65").Value = 9332.556
$ws.Range("I65").Value = 7999.25
$ws.Range("K65").Value = 39996.25
$ws.Range("M65").Value = -36876.25
$ws.Range("H107").Value = 3045.6667
$ws.Range("I107").Value = 4383.1665
$ws.Range("K107").Value = 4383.1665
$ws.Range("M107").Value = -2463.1665
$ws.Range("H113").Value = 2606.7273
$ws.Range("I113").Value = 1794.8
$ws.Range("K113").Value = 1794.8
$ws.Range("M113").Value = 375.2
$ws.Range("H132").Value = 2533.3044
$ws.Range("I132").Value = 2298.762
$ws.Range("J132").Value = 4996
$ws.Range("K132").Value = 6896.286
$ws.Range("L132").Value = 14988
$ws.Range("M132").Value = -4366.286
$ws.Range("N132").Value = -20048

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1074
$ws.Range("I68").Value = 652.5
$ws.Range("J68").Value = 2198
$ws.Range("K68").Value = 1957.5
$ws.Range("L68").Value = 6594
$ws.Range("M68").Value = -1146.5
$ws.Range("N68").Value = -8216
$ws.Range("H71").Value = 1074
$ws.Range("I71").Value = 652.5
$ws.Range("J71").Value = 2198
$ws.Range("K71").Value = 5872.5
$ws.Range("L71").Value = 19782
$ws.Range("M71").Value = -1816.5
$ws.Range("N71").Value = -27894
$ws.Range("H80").Value = 4832.3335
$ws.Range("J80").Value = 5198.8
$ws.Range("L80").Value = 15596.4
$ws.Range("N80").Value = -17468.4
$ws.Range("H83").Value = 4832.3335
$ws.Range("J83").Value = 5198.8
$ws.Range("L83").Value = 46789.2
$ws.Range("N83").Value = -56149.2
$ws.Range("H92").Value = 237.8
$ws.Range("I92").Value = 214.28572
$ws.Range("K92").Value = 642.85716
$ws.Range("M92").Value = 605.14284
$ws.Range("H131").Value = 18409.3
$ws.Range("J131").Value = 36029
$ws.Range("L131").Value = 108087
$ws.Range("N131").Value = -118167

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1999.9286
$ws.Range("I102").Value = 1335
$ws.Range("K102").Value = 1335
$ws.Range("M102").Value = 287

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H53").Value = 9496
$ws.Range("I53").Value = 6996.3335
$ws.Range("K53").Value = 6996.3335
$ws.Range("M53").Value = -6478.3335
$ws.Range("H61").Value = 757.2857
$ws.Range("I61").Value = 576
$ws.Range("K61").Value = 576
$ws.Range("M61").Value = -374
$ws.Range("H113").Value = 757.2857
$ws.Range("I113").Value = 576
$ws.Range("K113").Value = 576
$ws.Range("M113").Value = 1594
$ws.Range("H132").Value = 6390.8184
$ws.Range("I132").Value = 3143.75
$ws.Range("K132").Value = 9431.25
$ws.Range("M132").Value = -6901.25
$ws.Range("H136").Value = 2022.4193
$ws.Range("I136").Value = 2039.4828
$ws.Range("K136").Value = 6118.4484
$ws.Range("M136").Value = -3568.4484

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 49749.75
$ws.Range("J31").Value = 49749.75
$ws.Range("L31").Value = 49749.75
$ws.Range("N31").Value = -50445.75
$ws.Range("H107").Value = 999
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 999
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 2997
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -6837
$ws.Range("H122").Value = 1383
$ws.Range("I122").Value = 1295.2222
$ws.Range("K122").Value = 3885.6666
$ws.Range("M122").Value = -1435.6666
$ws.Range("H132").Value = 3504.2
$ws.Range("I132").Value = 3419.7
$ws.Range("J132").Value = 4011.2
$ws.Range("K132").Value = 10259.1
$ws.Range("L132").Value = 12033.6
$ws.Range("M132").Value = -7729.099999999999
$ws.Range("N132").Value = -17093.6

